# Apply betexplorer scrape update (script run 21-12-2023 02:45):
#  - Rows 84/85 had their match data (columns F:V) swapped back to the
#    correct order (Mladost vs Radnicki Nis should be row 84, Vojvodina
#    vs Radnik should be row 85).
#  - Rows 90/91 had their match data (columns F:V) swapped back to the
#    correct order (Radnicki Nis vs Sp. Subotica should be row 90,
#    IMT Novi Beograd vs Crvena zvezda should be row 91).
#  - A new match row (Partizan vs Crvena zvezda) was appended as row 145.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($ws, $rowA, $rowB, $firstCol, $lastCol) {
    $rangeA = $ws.Range("$firstCol$rowA" + ":" + "$lastCol$rowA")
    $rangeB = $ws.Range("$firstCol$rowB" + ":" + "$lastCol$rowB")

    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2

    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# Swap the two pairs of rows back to the correct match order.
Swap-RowRange $ws 84 85 "F" "V"
Swap-RowRange $ws 90 91 "F" "V"

# Append the new match as row 145. Copy formatting from the row above
# first (columns A and E carry non-default styles: bold/bordered index
# cell and a datetime number format) so the new cells reuse the exact
# same style records instead of Excel creating new ones.
$newRow = 145
$ws.Range("A144").Copy($ws.Range("A145"))
$ws.Range("E144").Copy($ws.Range("E145"))

$ws.Cells.Item($newRow, 1).Value2 = 144
$ws.Cells.Item($newRow, 2).Value2 = "serbia"
$ws.Cells.Item($newRow, 3).Value2 = "super-liga"
$ws.Cells.Item($newRow, 4).Value2 = "2023-2024"
$ws.Cells.Item($newRow, 5).Value2 = 45280.75
$ws.Cells.Item($newRow, 6).Value2 = "Partizan"
$ws.Cells.Item($newRow, 7).Value2 = 2
$ws.Cells.Item($newRow, 8).Value2 = "Crvena zvezda"
$ws.Cells.Item($newRow, 9).Value2 = 1
$ws.Cells.Item($newRow, 10).Value2 = 3.86
$ws.Cells.Item($newRow, 11).Value2 = "26/09/2023 07:12"
$ws.Cells.Item($newRow, 12).Value2 = 3.92
$ws.Cells.Item($newRow, 13).Value2 = "20/12/2023 17:59"
$ws.Cells.Item($newRow, 14).Value2 = 3.53
$ws.Cells.Item($newRow, 15).Value2 = "26/09/2023 07:12"
$ws.Cells.Item($newRow, 16).Value2 = 3.33
$ws.Cells.Item($newRow, 17).Value2 = "20/12/2023 17:59"
$ws.Cells.Item($newRow, 18).Value2 = 1.75
$ws.Cells.Item($newRow, 19).Value2 = "26/09/2023 07:12"
$ws.Cells.Item($newRow, 20).Value2 = 1.93
$ws.Cells.Item($newRow, 21).Value2 = "20/12/2023 17:57"
$ws.Cells.Item($newRow, 22).Value2 = "https://www.betexplorer.com/football/serbia/super-liga/partizan-crvena-zvezda/WbJ4er0B/"
